$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996111930376.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961119544177.csv"
$ws1.Range("B4").Value = "go_stims-16509961119544177.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961119704132.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509961133784156.csv"
$ws2.Range("B3").Value = "ZB-match_9-16509961127543786.csv"
$ws2.Range("B4").Value = "ZB-match_0-1650996112562412.csv"
$ws2.Range("B5").Value = "OB-16509961131944103.csv"
$ws2.Range("B6").Value = "ZB-match_4-16509961124183767.csv"
$ws2.Range("B7").Value = "TB-1650996113722417.csv"
$ws2.Range("B8").Value = "OB-16509961128263743.csv"
$ws2.Range("B9").Value = "OB-16509961131624124.csv"
$ws2.Range("B10").Value = "TB-1650996113738425.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996113770411.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996113754445.csv"
$ws4.Range("B4").Value = "MM_stims-16509961137863753.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996113770411.csv"
$ws4.Range("B6").Value = "MM_stims-16509961138023832.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961137863753.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961138184118.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961138344178.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961138023832.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961138503747.csv"

# --- Rename sheets (after content changes so names used above refer to positions, not old names) ---
$ws1.Name = "GNG_TO-16509961119704132"
$ws2.Name = "NB_TO-1650996113754445"
$ws3.Name = "RS_TO-1650996113754445"
$ws4.Name = "TOL_TO-16509961138023832"
$ws5.Name = "vSAT_TO-16509961138664203"
